# BM Table.xlsx update
#  - gssa.py folded into BMsolveGSSA.py (no workbook-visible effect).
#  - VFI renamed to VFI_11: the VFI sheet / Table "VFI" column are filled
#    in with the (11 grid point) results.
#  - Results for the 21- and 31-grid-point VFI runs are not carried in
#    this workbook state (only the numbers that belong here are written).
#
# NOTE: None of this touches the selections more than necessary - cells
# that already carried a style in the template (Table!E3/E4/E6/E7, and
# the Euler-error row 9 default format) simply get a Value/Formula so the
# existing style is kept; only genuinely new/re-typed cells get an
# explicit NumberFormat so the workbook picks up (or creates) the right
# style index.

$wb = $excel.ActiveWorkbook

$wsTable = $wb.Worksheets.Item("Table")
$wsEX    = $wb.Worksheets.Item("EX")
$wsLIN   = $wb.Worksheets.Item("LIN")
$wsGSSA  = $wb.Worksheets.Item("GSSA")
$wsVFI   = $wb.Worksheets.Item("VFI")

# ---------------------------------------------------------------------
# 1. Populate the VFI sheet (rows 3-10, columns B:E) with the VFI_11
#    (11 grid point) results. These are percentages, same number format
#    as the other "moments" sheets (EX / LIN / GSSA).
# ---------------------------------------------------------------------
$wsVFI.Range("B3").Value = 0.0277672
$wsVFI.Range("C3").Value = 0.39172000000000001
$wsVFI.Range("D3").Value = 0.0597436
$wsVFI.Range("E3").Value = 0.50685500000000006

$wsVFI.Range("B4").Value = 0.00943917
$wsVFI.Range("C4").Value = 0.13103600000000001
$wsVFI.Range("D4").Value = 0.00725202
$wsVFI.Range("E4").Value = 0.00908927

$wsVFI.Range("B5").Value = 0.0152861
$wsVFI.Range("C5").Value = 0.212204
$wsVFI.Range("D5").Value = 0.0117442
$wsVFI.Range("E5").Value = 0.0147195

$wsVFI.Range("B6").Value = 0.000234266
$wsVFI.Range("C6").Value = 0.015229
$wsVFI.Range("D6").Value = 0.00394682
$wsVFI.Range("E6").Value = 0.00494672

$wsVFI.Range("B7").Value = 0.18878300000000001
$wsVFI.Range("C7").Value = 2.62071999999999994
$wsVFI.Range("D7").Value = 0.14504
$wsVFI.Range("E7").Value = 0.181785

$wsVFI.Range("B8").Value = 0.0154199
$wsVFI.Range("C8").Value = 0.212559
$wsVFI.Range("D8").Value = 0.0140787
$wsVFI.Range("E8").Value = 0.0176484

$wsVFI.Range("B9").Value = 0.0279055
$wsVFI.Range("C9").Value = 0.39269399999999999
$wsVFI.Range("D9").Value = 0.00904901
$wsVFI.Range("E9").Value = 0.0118252

$wsVFI.Range("B10").Value = 0.00452137
$wsVFI.Range("C10").Value = 0.0648523
$wsVFI.Range("D10").Value = 0.000875539
$wsVFI.Range("E10").Value = 0.00112542

$wsVFI.Range("B3:E10").NumberFormat = "0.00%"

# ---------------------------------------------------------------------
# 2. Wire the Table sheet's VFI column (E) up to the new VFI numbers.
#    These cells already carry the shared "percent" style (s=2) in the
#    template, so only the formula needs to be written.
# ---------------------------------------------------------------------
$wsTable.Range("E3").Formula = "=VFI!B3"
$wsTable.Range("E4").Formula = "=VFI!C3"
$wsTable.Range("E6").Formula = "=VFI!D3"
$wsTable.Range("E7").Formula = "=VFI!E3"

# Euler error row (row 9) - VFI column. Row 9 already carries a custom
# row format, so the new cell just needs its value.
$wsTable.Range("E9").Value = 0.0010462799999999999

# ---------------------------------------------------------------------
# 3. Solve / Simulate timings (rows 11-12) - add the VFI column and
#    restyle the "Simulate" row (12) with a 2-decimal number format.
# ---------------------------------------------------------------------
$wsTable.Range("D11").Value = 0.2924845068967401
$wsTable.Range("D11").NumberFormat = "0.0000"
$wsTable.Range("E11").Value = 0.83543701622406852
$wsTable.Range("E11").NumberFormat = "0.0000"

$wsTable.Range("B12").Value = 31.537829729678734
$wsTable.Range("C12").Value = 62.380262520635483
$wsTable.Range("D12").Value = 122.76904415020975
$wsTable.Range("E12").Value = 41.49669063586263462
$wsTable.Range("B12:E12").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# 4. Selection bookkeeping left behind by the author's last save.
# ---------------------------------------------------------------------
$wsVFI.Activate()
$wsVFI.Range("E17").Select() | Out-Null

$wsTable.Activate()
$wsTable.Range("D17").Select() | Out-Null
